$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 389.7088409463599
$ws.Range("C2").Value = 30.49981016068242
$ws.Range("D2").Value = 470.8787783211095
